# Apply updated cryptocurrency price/volume data (scheduled data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "22.411.64"
$ws.Range("E2").Value = "  +0.00%  "

# Row 3
$ws.Range("D3").Value = "1.567.58"
$ws.Range("E3").Value = "  +0.06%  "

# Row 5
$ws.Range("E5").Value = "  +0.00%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "286.97"
$ws.Range("E6").Value = "  +0.78%  "

# Row 7
$ws.Range("E7").Value = "  +1.78%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "46.93"
$ws.Range("E8").Value = "  -3.37%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3299"
$ws.Range("E9").Value = "  -0.87%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.154"
$ws.Range("E10").Value = "  +2.72%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07447"
$ws.Range("E11").Value = "  +0.70%  "

# Row 12
$ws.Range("E12").Value = "  -0.02%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.61"
$ws.Range("E13").Value = "  -0.57%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.872"
$ws.Range("E14").Value = "  -1.29%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.839"
$ws.Range("E15").Value = "  -0.85%  "

# Row 16
$ws.Range("D16").Value = "1.566.50"
$ws.Range("E16").Value = "  +0.02%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001107"
$ws.Range("E17").Value = "  +0.36%  "

# Row 18
$ws.Range("E18").Value = "  -0.05%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "86.56"
$ws.Range("E19").Value = "  -1.67%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9998"
$ws.Range("E20").Value = "  -0.12%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.353"
$ws.Range("E21").Value = "  +0.07%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.39"
$ws.Range("E22").Value = "  +1.36%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.88"
$ws.Range("E23").Value = "  -0.97%  "

# Row 24
$ws.Range("D24").Value = "22.395.37"
$ws.Range("E24").Value = "  -0.07%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.333"
$ws.Range("E25").Value = "  -2.16%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.587"
$ws.Range("E26").Value = "  +1.61%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "151.39"
$ws.Range("E27").Value = "  +0.95%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.46"
$ws.Range("E28").Value = "  +0.61%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.938"
$ws.Range("E29").Value = "  -1.14%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "124.20"
$ws.Range("E30").Value = "  +0.38%  "

# Row 31
$ws.Range("D31").Value = "1.741.88"
$ws.Range("E31").Value = "  -0.04%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.059"
$ws.Range("E32").Value = "  +0.56%  "

# Row 33
$ws.Range("B33").Value = "WEMIXTOKEN"
$ws.Range("C33").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.975"
$ws.Range("E33").Value = "  -0.78%  "

# Row 34
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.021"
$ws.Range("E34").Value = "  -1.28%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.767"
$ws.Range("E35").Value = "  -0.26%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08275"
$ws.Range("E36").Value = "  +0.02%  "

# Row 37
$ws.Range("E37").Value = "  +0.08%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06363"
$ws.Range("E38").Value = "  -0.82%  "

# Row 39
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.292"
$ws.Range("E39").Value = "  +0.93%  "

# Row 40
$ws.Range("B40").Value = "Algorand"
$ws.Range("C40").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2199"
$ws.Range("E40").Value = "  -1.54%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.271"
$ws.Range("E41").Value = "  -1.74%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.29"
$ws.Range("E42").Value = "  +1.30%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6148"
$ws.Range("E43").Value = "  -1.50%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").Value = "  -0.06%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.83"
$ws.Range("E45").Value = "  +0.63%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5978"
$ws.Range("E46").Value = "  -0.85%  "

# Row 47
$ws.Range("E47").Value = "  +0.23%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.029"
$ws.Range("E48").Value = "  -0.10%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "123.83"
$ws.Range("E49").Value = "  -0.01%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.190"
$ws.Range("E50").Value = "  -1.88%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07166"
$ws.Range("E51").Value = "  -0.48%  "
